$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 13 and 14 follow the same visual pattern as the existing data
# rows (A/B/C/J/K use cell style index "1", the default body-text style).
# Copy that formatting from existing cells that already carry style "1"
# onto the new cells before setting their values, so no brand-new style
# gets created in styles.xml.

$ws.Range("A5:C5").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)

$ws.Range("J10").Copy()
$ws.Range("J13:J14").PasteSpecial(-4122)

$ws.Range("K5").Copy()
$ws.Range("K13:K14").PasteSpecial(-4122)

$ws.Range("A5:C5").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)

$ws.Range("A13").Value = "coupon_id"
$ws.Range("B13").Value = "No"
$ws.Range("C13").Value = "string"
$ws.Range("J13").Value = "UUID"
$ws.Range("K13").Value = "f602900a-4d93-4c44-a4a0-82e03a93d769"

$ws.Range("A14").Value = "campaign_id"
$ws.Range("B14").Value = "No"
$ws.Range("C14").Value = "string"
$ws.Range("J14").Value = "UUID"
$ws.Range("K14").Value = "f602900a-4d93-4c44-a4a0-82e03a93d769"
